$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated evaluation values (ifoCAST full series evaluation)
$ws.Range("B2").Value = 0.3341370637833943
$ws.Range("C2").Value = 0.5466446010724924
$ws.Range("D2").Value = 0.3940544821317493
$ws.Range("E2").Value = 0.6277375901853809
$ws.Range("F2").Value = 0.5468263669620367
$ws.Range("G2").Value = 18

$ws.Range("B3").Value = 0.1764647876375574
$ws.Range("C3").Value = 0.4397909425088788
$ws.Range("D3").Value = 0.3030744539125594
$ws.Range("E3").Value = 0.5505219831328804
$ws.Range("F3").Value = 0.5375226015493471
$ws.Range("G3").Value = 17

$ws.Range("B4").Value = 0.2016701448718926
$ws.Range("C4").Value = 0.3903364767990596
$ws.Range("D4").Value = 0.2840208244252586
$ws.Range("E4").Value = 0.5329360415896627
$ws.Range("F4").Value = 0.5094833745067471
$ws.Range("G4").Value = 16

$ws.Range("B5").Value = 0.3707750704723587
$ws.Range("C5").Value = 0.3719438412346756
$ws.Range("D5").Value = 0.1930669668093154
$ws.Range("E5").Value = 0.4393938629627357
$ws.Range("F5").Value = 0.2440568155285322
$ws.Range("G5").Value = 15

$ws.Range("B6").Value = 0.3163821704670755
$ws.Range("C6").Value = 0.3163821704670755
$ws.Range("D6").Value = 0.1380687648287419
$ws.Range("E6").Value = 0.3715760552413757
$ws.Range("F6").Value = 0.2022175558364308
$ws.Range("G6").Value = 14

$ws.Range("B7").Value = 0.3236369919384329
$ws.Range("C7").Value = 0.3236369919384329
$ws.Range("D7").Value = 0.1436765385716536
$ws.Range("E7").Value = 0.3790468817595702
$ws.Range("F7").Value = 0.2053783638939465
$ws.Range("G7").Value = 13

$ws.Range("B8").Value = 0.3336901046889824
$ws.Range("C8").Value = 0.3404581697164374
$ws.Range("D8").Value = 0.152137479526204
$ws.Range("E8").Value = 0.3900480477148989
$ws.Range("F8").Value = 0.2109417676443862
$ws.Range("G8").Value = 12

$ws.Range("B9").Value = 0.3630110268050817
$ws.Range("C9").Value = 0.3630110268050817
$ws.Range("D9").Value = 0.1611799638335552
$ws.Range("E9").Value = 0.4014722454087645
$ws.Range("F9").Value = 0.1798423033566436
$ws.Range("G9").Value = 11

$ws.Range("B10").Value = 0.3449376891861548
$ws.Range("C10").Value = 0.3449376891861548
$ws.Range("D10").Value = 0.1503712668170814
$ws.Range("E10").Value = 0.3877773418046514
$ws.Range("F10").Value = 0.1867537219501101
$ws.Range("G10").Value = 10

$ws.Range("B11").Value = 0.3563398907253821
$ws.Range("C11").Value = 0.3563398907253821
$ws.Range("D11").Value = 0.1562081803353234
$ws.Range("E11").Value = 0.3952318058245355
$ws.Range("F11").Value = 0.1813389655859693
$ws.Range("G11").Value = 9
